$d = $word.ActiveDocument

$p1Range = $d.Paragraphs(1).Range
$p1Xml = '<w:p w14:paraId="3E6B5F9E" w14:textId="425AFA5F" w:rsidR="00010884" w:rsidRDefault="00E01905" w:rsidP="006601EA" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>FSC</w:t></w:r><w:r w:rsidR="006601EA"><w:t xml:space="preserve"> QGIS </w:t></w:r><w:r><w:t>plugin help</w:t></w:r></w:p>'
$p1Range.InsertXML($p1Xml)

$p2Range = $d.Paragraphs(2).Range
$p2Xml = '<w:p w14:paraId="383FFC74" w14:textId="756AB7A2" w:rsidR="006601EA" w:rsidRDefault="00E01905" w:rsidP="006601EA" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:r><w:t>Help</w:t></w:r><w:r w:rsidR="00CF4B7C"><w:t xml:space="preserve"> for this </w:t></w:r><w:r><w:t>plugin</w:t></w:r><w:r w:rsidR="00CF4B7C"><w:t xml:space="preserve"> is available here</w:t></w:r><w:r w:rsidR="00A92E56"><w:t xml:space="preserve">: </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:instrText xml:space="preserve"> HYPERLINK "</w:instrText></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:instrText>https://www.fscbiodiversity.uk/qgisplugin</w:instrText></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:instrText xml:space="preserve">" </w:instrText></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://www.fscbiodiversity.uk/qgisplugin</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p>'
$p2Range.InsertXML($p2Xml)
